# Add test indexed_non-1 price
# Adds a new worksheet "TwoxTwowOTax_IndPrice_Nest" at the end of the workbook
# containing the benchmark / counterfactual results table for the indexed
# non-unity price test case, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# New sheet goes after the current last sheet, so it lands at the end of the
# tab strip (sheetId/rId are assigned sequentially by the engine).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TwoxTwowOTax_IndPrice_Nest"

# ---- Data rows: Column A / B labels + five numeric columns (C-G) ----
$rows = @(
    @("'x",       "'L",  50,  52.440442414743885,  52.440442408506513,  52.440442406059617,  51.056545362143005),
    @("'x",       "'K",  50,  47.673129456610248,  47.673129462280592,  47.673129464505024,  48.965318398797884),
    @("'y",       "'L",  20,  21.177057061457788,  21.177057058435174,  21.177057057249421,  20.508207428582388),
    @("'y",       "'K",  30,  28.87780507693768,   28.877805079685501,  28.877805080763462,  29.502316482825769),
    @("'Y.L",     "'x",  1,   1.0488088482223459,  1.0488088481702902,  1.0488088481593993,  1.256776502947158),
    @("'Y.L",     "'y",  1,   1.0388601182966062,  1.0388601182538628,  1.0388601182505186,  0.6257657257513739),
    @("'U.L",     "'_",  1,   1.0454820636283344,  1.0454820635787618,  1.0454820635708382,  1.0931727477635904),
    @("'PC.L",    "'x",  1,   0.99795751000197275, 1,                   1.0488088481386344,  1.0211309078784696),
    @("'PC.L",    "'y",  1,   1.0075145327263733,  1.0095765827768823,  1.0588528528793757,  1.0254103720418588),
    @("'PF.L",    "'L",  1,   0.95151515128525777, 0.95346258924559502, 1,                   1),
    @("'PF.L",    "'K",  1,   1.0466666666627262,  1.0488088481701119,  1.0999999998973025,  1.0427083297266266),
    @("'PU.L",    "'_",  1,   1.0011330687136686,  1.0031820580257127,  1.0521462187479673,  0.97829409400508083),
    @("'DU.L",    "'x",  100, 100.31820580330329,  100.31820580257127,  100.318205802237,    114.96595624939472),
    @("'DU.L",    "'y",  50,  49.683306602247974,  49.68330660297304,   49.683306603304153,  28.6215388689448),
    @("'SU.L",    "'_",  150, 150,                 150,                 150,                 150),
    @("'SY.L(i)", "'x",  100, 100,                 100,                 100,                 100),
    @("'SY.L(i)", "'y",  50,  50,                  50,                  50,                  50),
    @("'RA.L",    "'_",  150, 157,                 157.32132722551978,  164.99999999178419,  160.41666637813012)
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $r = $r + 1
}

# ---- Header row (typed last, G/E/F order) ----
$ws.Range("C1").Value = "'benchmark"
$ws.Range("D1").Value = "'RA=157"
$ws.Range("G1").Value = "'Pr.x=2"
$ws.Range("E1").Value = "'PC.x=1"
$ws.Range("F1").Value = "'PF.l=1"

# Make the new sheet the active tab (matches the workbook's activeTab bump).
$ws.Activate()
